$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.618.54'
$ws.Range('E2').Value = '  -1.82%  '

$ws.Range('D3').Value = '3.151.88'
$ws.Range('E3').Value = '  -4.68%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.02%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.68'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.71%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').Value = '3.148.18'
$ws.Range('E8').Value = '  -4.81%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.507'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.46%  '

$ws.Range('E10').Value = '  -5.62%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.26'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.12%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.453'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.80%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000234'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.64%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.80%  '

$ws.Range('D15').Value = '3.641.49'
$ws.Range('E15').Value = '  -5.61%  '

$ws.Range('E16').Value = '  -1.78%  '

$ws.Range('D17').Value = '3.165.79'
$ws.Range('E17').Value = '  -4.23%  '

$ws.Range('D18').Value = '62.560.36'
$ws.Range('E18').Value = '  -2.09%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.60%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '453.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.64%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.56%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.699'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.59%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.58'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.14%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.55%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.43'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.28%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.43%  '

$ws.Range('B27').Value = 'FirstDigitalUSD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.14%  '

$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.68'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.27%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.68'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.52%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.37%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.01'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.92%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.16'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.61%  '

$ws.Range('E33').Value = '  -3.47%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.37'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.09%  '

$ws.Range('E35').Value = '  -7.24%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.80'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.47%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.10'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.23%  '

$ws.Range('D38').Value = '0.0₃0702'
$ws.Range('E38').Value = '  -5.60%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0385'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.01%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '401.43'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.38%  '

$ws.Range('E41').Value = '  -3.21%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.01'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.06%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.110'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.91%  '

$ws.Range('D44').Value = '2.766.83'
$ws.Range('E44').Value = '  -10.12%  '

$ws.Range('E45').Value = '  -6.05%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.11'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.92%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.38'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.01%  '

$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.18'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.45%  '

$ws.Range('B50').Value = 'Arweave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.59'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.15%  '

$ws.Range('E51').Value = '  -3.61%  '
